# Profile updates week 16
#
# Flips the "active" indicator (column C) for three profile rows and
# updates their cell style to the one already used by other "flipped"
# rows in the sheet (cellXfs index 3, as seen e.g. on C89):
#   - Row 10 (Manas / manas_rishav):          1 -> 0
#   - Row 82 (Kiran S / kiran_subramoni1):    0 -> 1
#   - Row 88 (Rishabh Gupta / rishabh_gupta4):0 -> 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C89 already carries the target cell style; copy just its formatting
# onto each target cell, then set the new flag value.
$styleSource = $ws.Range("C89")

$styleSource.Copy()
$ws.Range("C10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C10").Value = 0

$styleSource.Copy()
$ws.Range("C82").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C82").Value = 1

$styleSource.Copy()
$ws.Range("C88").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C88").Value = 1

$excel.CutCopyMode = $false
